$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 data: December entry
$ws.Range("A4").Value = 2024
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = "December"
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 23
$ws.Range("F4").Formula = "=E4-D4+1"

# Update selection to E8 as in the diff
$ws.Range("E8").Select()
